$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<name>_old" -> "<name>_FV2304" (cols A-J),
#    "<name>_new" -> "<name>_FV2310" (cols L-U). Column K ("diff") is kept.
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2304"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2310"
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U94 into an Excel Table (ListObject), keeping the header row's
#    existing look (bold font, grey fill, centered/wrapped, thin border)
#    instead of letting Excel capture it into a new header-row dxf.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")

# stash a copy of the still-formatted header style on a scratch cell
$ws.Range("A1").Copy($ws.Range("W1"))

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# restore the original header formatting
$ws.Range("W1").Copy()
$headerRange.PasteSpecial(-4122)
$ws.Range("W1").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
